# Adds the "30. 11. 2021" wave of data to both worksheets:
#   - "data"   : new column AK (+ two small corrections to the previous AJ column)
#   - "pocetR" : new column AJ (+ small corrections to the previous AI column)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data" (percentages)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# New header cell AK1, copying the header formatting used by the rest of row 1.
$ws1.Range("AJ1").Copy() | Out-Null
$ws1.Range("AK1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws1.Range("AK1").Value = "30. 11. 2021"

# New values for column AK, rows 2-45.
$ak1Data = @(
    0.22, 0.13, 0.46, 0.31, 0.15, 0.2,  0.31, 0.23, 0.2,  0.2,
    0.27, 0.44, 0.2,  0.21, 0.25, 0.19, 0.23, 0.28, 0.19, 0.16,
    0.18, 0.2,  0.38, 0.44, 0.13, 0.09, 0.14, 0.18, 0.11, 0.12,
    0.14, 0.16, 0.22, 0.11, 0.16, 0.13, 0.07000000000000001, 0.21,
    0.15, 0.12, 0.1,  0.07000000000000001, 0.18, 0.28
)
for ($i = 0; $i -lt $ak1Data.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 37).Value = $ak1Data[$i]
}

# Small corrections made to the previous (AJ) wave at the same time.
$ws1.Range("AJ20").Value = 0.21
$ws1.Range("AJ40").Value = 0.13

# Footer note row - bump the "aktualizace" date.
$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 8. 12. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" (sample sizes)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# New header cell AJ1, copying the header formatting used by the rest of row 1.
$ws2.Range("AI1").Copy() | Out-Null
$ws2.Range("AJ1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws2.Range("AJ1").Value = "30. 11. 2021"

# New values for column AJ, rows 2-23.
$aj2Data = @(
    1790, 192, 352, 1246, 865, 159, 512, 254, 829, 144,
    113,  704, 829, 608,  353, 204, 633, 585, 266, 521,
    320,  177
)
for ($i = 0; $i -lt $aj2Data.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 36).Value = $aj2Data[$i]
}

# Small corrections made to the previous (AI) wave at the same time.
$ws2.Range("AI3").Value = 166
$ws2.Range("AI4").Value = 366
$ws2.Range("AI17").Value = 184
$ws2.Range("AI18").Value = 621
$ws2.Range("AI19").Value = 537
$ws2.Range("AI20").Value = 252

# Footer note row - bump the "aktualizace" date.
$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 8. 12. 2021"

# Trailing empty placeholder cell that mirrors the rest of the (blank) footer row.
$ws2.Range("AJ24").NumberFormat = "General"
